$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-18 are bumped from 45179 to 45180 (one day later)
$ws.Range("C2:C18").Value = 45180
